$wb = $excel.ActiveWorkbook

# --- TestSuite: flip RunMode flags (enable LoginPageTest, disable TimesheetPageTest) ---
$ws1 = $wb.Worksheets.Item("TestSuite")
$ws1.Range("B2").Value = "Y"
$ws1.Range("B4").Value = "N"

# --- LoginPageTest: replace plaintext passwords with encrypted values ---
$ws2 = $wb.Worksheets.Item("LoginPageTest")
$ws2.Range("B2").Value = "gG+7Twxtcof2boCuiDPlzA=="
$ws2.Range("B3").Value = "/bbmmvb4w8JDpN0RgERT3w=="
$ws2.Columns.Item(2).ColumnWidth = 28.43

# --- AdminPageTest: replace plaintext password with encrypted value ---
$ws3 = $wb.Worksheets.Item("AdminPageTest")
$ws3.Range("B2").Value = "gG+7Twxtcof2boCuiDPlzA=="
$ws3.Columns.Item(2).ColumnWidth = 25.33

# --- TimesheetPageTest: replace plaintext password with encrypted value ---
$ws4 = $wb.Worksheets.Item("TimesheetPageTest")
$ws4.Range("B2").Value = "gG+7Twxtcof2boCuiDPlzA=="
$ws4.Columns.Item(2).ColumnWidth = 25.33

# --- Update the active selection on each sheet, ending with TestSuite so it
#     becomes the active/selected tab (matches tabSelected moving from
#     TimesheetPageTest to TestSuite) ---
$ws2.Range("B3").Select()
$ws3.Range("E23").Select()
$ws4.Range("D23").Select()
$ws1.Range("B23").Select()
